$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column B header from "value" to "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# The series now extends from row 53 to row 84 (31 new rows). Clone the
# existing date-column formatting (bold font, thin border, centered/top
# aligned, custom date number format) from A53 onto the new A54:A84 cells
# so the whole date column stays visually consistent.
$ws.Range("A53").Copy()
$ws.Range("A54:A84").PasteSpecial(-4122)

# Rewrite the full two-column series: dates (column A) and values (column B)
$ws.Cells.Item(2, 1).Value = 38398
$ws.Cells.Item(2, 2).Value = -0.7
$ws.Cells.Item(3, 1).Value = 38487
$ws.Cells.Item(3, 2).Value = 0.4
$ws.Cells.Item(4, 1).Value = 38579
$ws.Cells.Item(4, 2).Value = -0.2
$ws.Cells.Item(5, 1).Value = 38671
$ws.Cells.Item(5, 2).Value = 0.7
$ws.Cells.Item(6, 1).Value = 38763
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(7, 1).Value = 38852
$ws.Cells.Item(7, 2).Value = -0.4
$ws.Cells.Item(8, 1).Value = 38944
$ws.Cells.Item(8, 2).Value = 0.1
$ws.Cells.Item(9, 1).Value = 39036
$ws.Cells.Item(9, 2).Value = -1.6
$ws.Cells.Item(10, 1).Value = 39128
$ws.Cells.Item(10, 2).Value = 2
$ws.Cells.Item(11, 1).Value = 39217
$ws.Cells.Item(11, 2).Value = -1.1
$ws.Cells.Item(12, 1).Value = 39309
$ws.Cells.Item(12, 2).Value = 0.4
$ws.Cells.Item(13, 1).Value = 39401
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(14, 1).Value = 39493
$ws.Cells.Item(14, 2).Value = 0.7
$ws.Cells.Item(15, 1).Value = 39583
$ws.Cells.Item(15, 2).Value = -0.4
$ws.Cells.Item(16, 1).Value = 39675
$ws.Cells.Item(16, 2).Value = 0.9
$ws.Cells.Item(17, 1).Value = 39767
$ws.Cells.Item(17, 2).Value = 0.6
$ws.Cells.Item(18, 1).Value = 39859
$ws.Cells.Item(18, 2).Value = -0.5
$ws.Cells.Item(19, 1).Value = 39948
$ws.Cells.Item(19, 2).Value = -2
$ws.Cells.Item(20, 1).Value = 40040
$ws.Cells.Item(20, 2).Value = 1.5
$ws.Cells.Item(21, 1).Value = 40132
$ws.Cells.Item(21, 2).Value = -1.1
$ws.Cells.Item(22, 1).Value = 40224
$ws.Cells.Item(22, 2).Value = 1.9
$ws.Cells.Item(23, 1).Value = 40313
$ws.Cells.Item(23, 2).Value = 0.1
$ws.Cells.Item(24, 1).Value = 40405
$ws.Cells.Item(24, 2).Value = -0.3
$ws.Cells.Item(25, 1).Value = 40497
$ws.Cells.Item(25, 2).Value = -0.5
$ws.Cells.Item(26, 1).Value = 40589
$ws.Cells.Item(26, 2).Value = -0.4
$ws.Cells.Item(27, 1).Value = 40678
$ws.Cells.Item(27, 2).Value = 0.3
$ws.Cells.Item(28, 1).Value = 40770
$ws.Cells.Item(28, 2).Value = -0.4
$ws.Cells.Item(29, 1).Value = 40862
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(30, 1).Value = 40954
$ws.Cells.Item(30, 2).Value = -0.4
$ws.Cells.Item(31, 1).Value = 41044
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(32, 1).Value = 41136
$ws.Cells.Item(32, 2).Value = -0.3
$ws.Cells.Item(33, 1).Value = 41228
$ws.Cells.Item(33, 2).Value = 0.4
$ws.Cells.Item(34, 1).Value = 41320
$ws.Cells.Item(34, 2).Value = -0.1
$ws.Cells.Item(35, 1).Value = 41409
$ws.Cells.Item(35, 2).Value = -0.3
$ws.Cells.Item(36, 1).Value = 41501
$ws.Cells.Item(36, 2).Value = 0.2
$ws.Cells.Item(37, 1).Value = 41593
$ws.Cells.Item(37, 2).Value = -0.2
$ws.Cells.Item(38, 1).Value = 41685
$ws.Cells.Item(38, 2).Value = 0.7000000000000001
$ws.Cells.Item(39, 1).Value = 41774
$ws.Cells.Item(39, 2).Value = 0.1
$ws.Cells.Item(40, 1).Value = 41866
$ws.Cells.Item(40, 2).Value = -0.5
$ws.Cells.Item(41, 1).Value = 41958
$ws.Cells.Item(41, 2).Value = 0.4
$ws.Cells.Item(42, 1).Value = 42050
$ws.Cells.Item(42, 2).Value = -0.3
$ws.Cells.Item(43, 1).Value = 42139
$ws.Cells.Item(43, 2).Value = -0.3
$ws.Cells.Item(44, 1).Value = 42231
$ws.Cells.Item(44, 2).Value = 0.2
$ws.Cells.Item(45, 1).Value = 42323
$ws.Cells.Item(45, 2).Value = 0.1
$ws.Cells.Item(46, 1).Value = 42415
$ws.Cells.Item(46, 2).Value = 0.1
$ws.Cells.Item(47, 1).Value = 42505
$ws.Cells.Item(47, 2).Value = -0.2
$ws.Cells.Item(48, 1).Value = 42597
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(49, 1).Value = 42689
$ws.Cells.Item(49, 2).Value = 0.4
$ws.Cells.Item(50, 1).Value = 42781
$ws.Cells.Item(50, 2).Value = -0.4
$ws.Cells.Item(51, 1).Value = 42870
$ws.Cells.Item(51, 2).Value = 0.2
$ws.Cells.Item(52, 1).Value = 42962
$ws.Cells.Item(52, 2).Value = 0.4
$ws.Cells.Item(53, 1).Value = 43054
$ws.Cells.Item(53, 2).Value = 0
$ws.Cells.Item(54, 1).Value = 43146
$ws.Cells.Item(54, 2).Value = -0.1
$ws.Cells.Item(55, 1).Value = 43235
$ws.Cells.Item(55, 2).Value = 0.4
$ws.Cells.Item(56, 1).Value = 43327
$ws.Cells.Item(56, 2).Value = 0.7
$ws.Cells.Item(57, 1).Value = 43419
$ws.Cells.Item(57, 2).Value = -0.6
$ws.Cells.Item(58, 1).Value = 43511
$ws.Cells.Item(58, 2).Value = -0.6
$ws.Cells.Item(59, 1).Value = 43600
$ws.Cells.Item(59, 2).Value = 0.3
$ws.Cells.Item(60, 1).Value = 43692
$ws.Cells.Item(60, 2).Value = -0.7
$ws.Cells.Item(61, 1).Value = 43784
$ws.Cells.Item(61, 2).Value = 0.6
$ws.Cells.Item(62, 1).Value = 43876
$ws.Cells.Item(62, 2).Value = 0.3
$ws.Cells.Item(63, 1).Value = 43966
$ws.Cells.Item(63, 2).Value = 0.3
$ws.Cells.Item(64, 1).Value = 44058
$ws.Cells.Item(64, 2).Value = -0.4440571223929872
$ws.Cells.Item(65, 1).Value = 44150
$ws.Cells.Item(65, 2).Value = -0.7255945204468831
$ws.Cells.Item(66, 1).Value = 44242
$ws.Cells.Item(66, 2).Value = -0.5292660609007143
$ws.Cells.Item(67, 1).Value = 44331
$ws.Cells.Item(67, 2).Value = -0.1550786956675604
$ws.Cells.Item(68, 1).Value = 44423
$ws.Cells.Item(68, 2).Value = -2.168330733759602
$ws.Cells.Item(69, 1).Value = 44515
$ws.Cells.Item(69, 2).Value = -0.03982694963614287
$ws.Cells.Item(70, 1).Value = 44607
$ws.Cells.Item(70, 2).Value = 0.2669401745841223
$ws.Cells.Item(71, 1).Value = 44696
$ws.Cells.Item(71, 2).Value = 0.03791487406588956
$ws.Cells.Item(72, 1).Value = 44788
$ws.Cells.Item(72, 2).Value = -0.04567208272808071
$ws.Cells.Item(73, 1).Value = 44880
$ws.Cells.Item(73, 2).Value = -0.5154625125417773
$ws.Cells.Item(74, 1).Value = 44972
$ws.Cells.Item(74, 2).Value = -0.1813602613933202
$ws.Cells.Item(75, 1).Value = 45061
$ws.Cells.Item(75, 2).Value = -0.01480819732384536
$ws.Cells.Item(76, 1).Value = 45153
$ws.Cells.Item(76, 2).Value = 0.02918400950819283
$ws.Cells.Item(77, 1).Value = 45245
$ws.Cells.Item(77, 2).Value = -0.03321544329283629
$ws.Cells.Item(78, 1).Value = 45337
$ws.Cells.Item(78, 2).Value = 0.00001303303454188581
$ws.Cells.Item(79, 1).Value = 45427
$ws.Cells.Item(79, 2).Value = -0.006125572440376981
$ws.Cells.Item(80, 1).Value = 45519
$ws.Cells.Item(80, 2).Value = 0.04879937325030748
$ws.Cells.Item(81, 1).Value = 45611
$ws.Cells.Item(81, 2).Value = 0.0477695913607396
$ws.Cells.Item(82, 1).Value = 45703
$ws.Cells.Item(82, 2).Value = 0.4714513528429705
$ws.Cells.Item(83, 1).Value = 45792
$ws.Cells.Item(83, 2).Value = -0.02605454389395597
$ws.Cells.Item(84, 1).Value = 45884
$ws.Cells.Item(84, 2).Value = 0.04549112474043772

Write-Host "Edit complete"
